# Fix typo in shared string "Doma}ska" -> "Domańska" and swap the
# first-name / last-name columns (A <-> B) for every data row, then
# move the active selection to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo wherever it appears on the sheet (column B originally).
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    if ($cell.Value2 -eq "Doma}ska") {
        $cell.Value2 = "Domańska"
    }
}

# Swap columns A and B for each data row (rows 1-6).
for ($r = 1; $r -le 6; $r++) {
    $a = $ws.Cells.Item($r, 1)
    $b = $ws.Cells.Item($r, 2)
    $tmp = $a.Value2
    $a.Value2 = $b.Value2
    $b.Value2 = $tmp
}

# Update the active selection to E3.
$ws.Range("E3").Select()
